# Apply the "check user in group / join group" edit to groups.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the three private-group invite links to show just the invite
# hash (the "HZesgX2L5zcpKvq0" style code) instead of the full
# https://t.me/joinchat/... URL.
$ws.Range("E2").Value = "HZesgX2L5zcpKvq0"
$ws.Range("E3").Value = "VTvg_eT6s7Rz-AIj"
$ws.Range("E4").Value = "RcGGtdG60NynCrJK"

# Update the group_id for the public "Test Tool" group.
$ws.Range("B5").Value = -1001159430667

# Turn the public group's link into a clickable hyperlink.
$ws.Hyperlinks.Add($ws.Range("E5"), "https://t.me/testInteractTool")

# Move the active selection to B5.
$ws.Range("B5").Select()

$wb.Save()
